$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Remove the stray "_GoBack" bookmark that sits between the ", " run and
# the "dst" run (near the author names), and merge the two runs into a
# single run with text ", dst".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Content.Find.Execute(", dst", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", dst", 2)

# --- Change 2 -------------------------------------------------------------
# Insert a new, empty paragraph (inheriting the paragraph's own formatting)
# right after the Indonesian abstract paragraph, before "Kata Kunci:".
$needleId = "Abstrak ditulis dalam bahasa Indonesia (kecuali artikel yang ditulis dalam bahasa Inggris) yang berisikan isu-isu pokok, tujuan penelitian, metode/pendekatan dan hasil penelitian. Abstrak ditulis dalam satu alinea, tidak lebih dari 250 kata. (Times New Roman 11, spasi 1)."
$d.Content.Find.Execute($needleId, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $needleId + "`r", 2)

# --- Change 3 -------------------------------------------------------------
# Same idea for the English abstract paragraph, before "Keywords:".
$needleEn = "Abstrak ditulis dalam bahasa Inggris yang berisikan isu-isu pokok, tujuan penelitian, metode/pendekatan dan hasil penelitian. Abstrak ditulis dalam satu alinea, tidak lebih dari 250 kata. (Times New Roman 11, spasi 1, dan cetak miring)."
$d.Content.Find.Execute($needleEn, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $needleEn + "`r", 2)

# --- Change 4 -------------------------------------------------------------
# Split "PENDAHULUAN" into "PENDA" + "HULUAN" with a collapsed "_GoBack"
# bookmark in between (mirrors the bookmark relocation caused by the
# author's last edit landing here).
$rng = $d.Content
$found = $rng.Find.Execute("PENDAHULUAN", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $rng.Start + 5
    $mid = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $mid)
}
